$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) figures.
# Column D values are forced as text (leading apostrophe) so Excel does not
# reinterpret dotted numeric strings (e.g. "600.95", "68.633.92") as floats
# and lose their original textual formatting/precision.

$ws.Range("D2").Value = "'68.633.92"
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").Value = "'3.842.74"
$ws.Range("E3").Value = '  +2.75%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'600.95"
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").Value = "'163.70"
$ws.Range("E6").Value = '  -2.29%  '
$ws.Range("D7").Value = "'3.841.43"
$ws.Range("E7").Value = '  +2.74%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  -2.11%  '
$ws.Range("E10").Value = '  -1.33%  '
$ws.Range("E11").Value = '  -0.38%  '
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("D13").Value = "'36.93"
$ws.Range("E13").Value = '  -3.23%  '
$ws.Range("E14").Value = '  -1.26%  '
$ws.Range("D15").Value = "'4.488.14"
$ws.Range("E15").Value = '  +2.77%  '
$ws.Range("D16").Value = "'3.878.95"
$ws.Range("E16").Value = '  +3.81%  '
$ws.Range("D17").Value = "'68.810.65"
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("D18").Value = "'7.57"
$ws.Range("E18").Value = '  +2.67%  '
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("D20").Value = "'17.16"
$ws.Range("E20").Value = '  -1.31%  '
$ws.Range("D21").Value = "'11.24"
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("D22").Value = "'486.58"
$ws.Range("E22").Value = '  -1.14%  '
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("E24").Value = '  +6.15%  '
$ws.Range("D25").Value = "'84.08"
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("E26").Value = '  -2.14%  '
$ws.Range("D27").Value = "'12.11"
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("D28").Value = "'10.00"
$ws.Range("E28").Value = '  -0.61%  '
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("E30").Value = '  -0.62%  '
$ws.Range("D31").Value = "'7.85"
$ws.Range("E31").Value = '  -3.98%  '
$ws.Range("D32").Value = "'3.994.36"
$ws.Range("E32").Value = '  +2.81%  '
$ws.Range("E33").Value = '  -3.85%  '
$ws.Range("E34").Value = '  +1.25%  '
$ws.Range("D35").Value = "'3.787.77"
$ws.Range("E35").Value = '  +3.12%  '
$ws.Range("E36").Value = '  -1.17%  '
$ws.Range("E37").Value = '  +1.83%  '
$ws.Range("E39").Value = '  -1.00%  '
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("E41").Value = '  -2.37%  '
$ws.Range("D42").Value = "'2.97"
$ws.Range("E42").Value = '  -2.05%  '
$ws.Range("D43").Value = "'432.19"
$ws.Range("E43").Value = '  +1.96%  '
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("D47").Value = "'8.41"
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("D48").Value = "'2.845.41"
$ws.Range("E48").Value = '  +2.32%  '
$ws.Range("D49").Value = "'142.71"
$ws.Range("E49").Value = '  +1.07%  '
$ws.Range("D50").Value = "'0.0358"
$ws.Range("E50").Value = '  +1.20%  '
$ws.Range("D51").Value = "'25.88"
$ws.Range("E51").Value = '  +13.28%  '
